$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Locate the "User story 6" body paragraph (currently split across two
#        runs around a _GoBack bookmark: "For the user, I w" + "ould like ...").
#        Find the paragraph whose range starts the merged text, using a
#        phrase that is unique to it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*add guideline while they are using the program*") {
        $target = $p
        break
    }
}

$mergedText = [char]0x201C + "New" + [char]0x201D + ", " + [char]0x201C + "Load map" + [char]0x201D
$bodyText = "For the user, I would like to add guideline while they are using the program to avoid confusion from the users. Priority is Low and estimate time is 1 day. For testing, press the " + $mergedText + " button will result in a short guideline."

$xml1 = "<w:p $wns w:rsidR=`"006265F4`" w:rsidRDefault=`"006265F4`" w:rsidP=`"001D7A14`"><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>$bodyText</w:t></w:r></w:p>"
$xml2 = "<w:p $wns><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>User story 7:</w:t></w:r></w:p>"
$gui = "For the user, I would like to change the GUI for better navigation. For testing, Users can now select buttons easier as all the buttons are now on the left side of the screen. "
$xml3 = "<w:p $wns><w:r><w:t xml:space=`"preserve`">$gui</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"

$target.Range.InsertXML($xml1 + $xml2 + $xml3) | Out-Null

# --- 2. Insert a new blank paragraph right after the developer story
#        paragraph (before the pre-existing blank paragraph).
$dev = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*keeping a distance between each other in the simulator*") {
        $dev = $p
        break
    }
}

$devXml = "<w:p $wns w:rsidR=`"00C57413`" w:rsidRDefault=`"00C57413`" w:rsidP=`"00C57413`"><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r w:rsidRPr=`"00C57413`"><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>For the developer, I would like to add some distance between vehicles to avoid collisions. In Australia rule, a vehicle should drive at least 2 second behind the vehicle in front. Priority is High and the estimate time is 1 day. For testing, see if the vehicle are keeping a distance between each other in the simulator.</w:t></w:r></w:p>"
$blankXml = "<w:p $wns><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr></w:p>"

$dev.Range.InsertXML($devXml + $blankXml) | Out-Null
